$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D21").Value = "[c++] Thread 사용"
$ws.Range("E21").Value = "https://ms-review.tistory.com/21"

$ws.Range("D36").Value = "손실함수(Loss function)"
$ws.Range("E36").Value = "http://dmqm.korea.ac.kr/activity/seminar/326"

$ws.Range("D37").Value = "[Paper Review] Knowledge Distillation 2021 ver."
$ws.Range("E37").Value = "http://dsba.korea.ac.kr/seminar/?uid=1789&mod=document&pageid=1"
